# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`r`n✅ Dólar paralelo: 68`r`n`r`nBinance`r`n✅ 1000 Bs = 3.71 = 14172.44 pesos`r`n✅ 14172.44 pesos = 3.69 = 963.31 Bs`r`n`r`nPromedio competencia`r`n✅ Tasa pesos: 20`r`n✅ Tasa Bs: 20`r`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas: update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 269.89
$ws2.Range("O10").Value = 3825
$ws2.Range("N12").Value = 3839.9
$ws2.Range("O12").Value = 261
